$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update comparison-matrix values (row 2 = Datasheet row, row 3 = Dataset Nutrition Label row, etc.)
# Diagonal cells (value 100) remain unchanged.

$ws.Range("B2").Value = 48.93617021276596
$ws.Range("C2").Value = 36.17021276595745
$ws.Range("D2").Value = 97.87234042553192
$ws.Range("E2").Value = 29.78723404255319

$ws.Range("A3").Value = 53.48837209302325
$ws.Range("C3").Value = 37.2093023255814
$ws.Range("D3").Value = 62.7906976744186
$ws.Range("E3").Value = 30.23255813953488

$ws.Range("A4").Value = 34.69387755102041
$ws.Range("B4").Value = 32.6530612244898
$ws.Range("D4").Value = 44.89795918367347
$ws.Range("E4").Value = 36.73469387755102

$ws.Range("A5").Value = 68.65671641791045
$ws.Range("B5").Value = 40.29850746268657
$ws.Range("C5").Value = 32.83582089552239
$ws.Range("E5").Value = 26.86567164179105

$ws.Range("A6").Value = 26.41509433962264
$ws.Range("B6").Value = 24.52830188679245
$ws.Range("C6").Value = 33.9622641509434
$ws.Range("D6").Value = 33.9622641509434
